$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header labels: "<Name>_old" -> "<Name>_FV2210" and
#    "<Name>_new" -> "<Name>_FV2304" (row 1, columns A:J and L:U; column K
#    stays "diff").
# ---------------------------------------------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2210")
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2304")
}

# ---------------------------------------------------------------------------
# 2. Turn the used range into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$rng = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row (pane split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
